$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert 4 rows before row 17 to push the totals rows (17,18) down to (21,22)
$ws.Rows("17:20").Insert()

# Fill in new row 16 - solder paste
$ws.Range("A16").Value = "solder paste"
$ws.Range("H16").Value = 12.82

# Fill in new row 17 - LED strip lights
$ws.Range("A17").Value = "LED strip lights"
$ws.Range("H17").Value = 73.96

# Update totals formulas (now at row 21)
$ws.Range("G21").Formula = "=SUM(G2:G15)"
$ws.Range("H21").Formula = "=SUM(H2:H17)"

# Update grand total formula (now at row 22)
$ws.Range("G22").Formula = "=G21+H21"

# Selection
$ws.Range("G22:H22").Select()
